# Generate Report for Handback
# Marks the two localization rows (zh-cn and de-de sheets) as handed back:
#  - Status text changes from "Ready for handoff" to "Handed back: in sync with en-US"
#    (this text is shared by the Overview sheet status columns too)
#  - The "Latest Target File" (I) / "Latest Handback File" (J) columns are filled in
#    with the handed-back .md file (with hyperlink) and the generated .xlf file name
#  - The de-de sheet also records a new "Latest Handback DateTime" (K)
#  - A few columns are widened so the new, longer text fits

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/756e7e0816735f27ec149b497ebe7151f6037179/e2e/02c6e247-5a6e-4c24-87b5-82a3b948c594.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/756e7e0816735f27ec149b497ebe7151f6037179/e2e/10826afa-e4d5-4f23-8910-2a2654cf7c00.md"

$mdName1 = "02c6e247-5a6e-4c24-87b5-82a3b948c594.md"
$mdName2 = "10826afa-e4d5-4f23-8910-2a2654cf7c00.md"

# Column width used by every widened column that previously held the
# (narrower) "17.2159881591797" width.
$wideWidth = 29.166666666666668
# Column width used for columns I and J on the language sheets (40 chars).
$fullWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# Overview sheet: status text + column widths
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText
$overview.Columns.Item(5).ColumnWidth = $wideWidth
$overview.Columns.Item(6).ColumnWidth = $wideWidth

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$zhcn.Range("I2").Value = $mdName1
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl1, "", "", $mdName1)
$zhcn.Range("J2").Value = "02c6e247-5a6e-4c24-87b5-82a3b948c594.42aa03bcf147d2a9e5f03d6c9ebb457ef8daa1da.zh-cn.xlf"

$zhcn.Range("I3").Value = $mdName2
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $mdUrl2, "", "", $mdName2)
$zhcn.Range("J3").Value = "10826afa-e4d5-4f23-8910-2a2654cf7c00.b5bc747da3cc2497c66fdafb0c20d2e0be701280.zh-cn.xlf"

# Latest Handback DateTime (K) for zh-cn was a placeholder value; fill it in
# for both data rows so the shared string is updated consistently.
$zhcn.Range("K2").Value = "2016-08-22 04:40:43"
$zhcn.Range("K3").Value = "2016-08-22 04:40:43"

$zhcn.Columns.Item(3).ColumnWidth = $wideWidth
$zhcn.Columns.Item(9).ColumnWidth = $fullWidth
$zhcn.Columns.Item(10).ColumnWidth = $fullWidth

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

$dede.Range("I2").Value = $mdName1
$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl1, "", "", $mdName1)
$dede.Range("J2").Value = "02c6e247-5a6e-4c24-87b5-82a3b948c594.42aa03bcf147d2a9e5f03d6c9ebb457ef8daa1da.de-de.xlf"
$dede.Range("K2").Value = "2016-08-22 04:40:50"

$dede.Range("I3").Value = $mdName2
$dede.Hyperlinks.Add($dede.Range("I3"), $mdUrl2, "", "", $mdName2)
$dede.Range("J3").Value = "10826afa-e4d5-4f23-8910-2a2654cf7c00.b5bc747da3cc2497c66fdafb0c20d2e0be701280.de-de.xlf"
$dede.Range("K3").Value = "2016-08-22 04:40:50"

$dede.Columns.Item(3).ColumnWidth = $wideWidth
$dede.Columns.Item(9).ColumnWidth = $fullWidth
$dede.Columns.Item(10).ColumnWidth = $fullWidth

Write-Host "Handback report generated"
